$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 28 with product data
$ws.Range("A28").Value = "AOQWST"
$ws.Range("B28").Value = "Chip Epson"
$ws.Range("C28").Value = "T04D1"
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 100000
$ws.Range("F28").Value = 28
$ws.Range("G28").Value = 2
$ws.Range("H28").Formula = "=(E28-D28)*G28"
$ws.Range("I28").Formula = "=D28*F28"
$ws.Range("J28").Value = 0
